# "Generate Report for Archive"
#
# The localization-status report text moved from "Ready for handoff" to
# "In Translation" (the same shared string is used for the per-language
# Status columns on the Overview sheet as well as on each language sheet),
# and the now-narrower text let the Status columns be shrunk.
#
# Note: reading a cell's value back out through `.Value` misbehaves when
# printed/compared directly in this host, and PowerShell's `-eq` coerces its
# right-hand side to the type of the left-hand side - cells holding the
# literal text "True"/"False" come back from `.Value`/`.Value2` as real
# Booleans, so comparing `$cell.Value2 -eq "Ready for handoff"` would coerce
# "Ready for handoff" (non-empty) to `$true` and spuriously match those
# cells. Keeping the literal string on the left of `-eq` avoids that.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    foreach ($cell in $ws.UsedRange.Cells) {
        if ("Ready for handoff" -eq $cell.Value2) {
            $cell.Value = "In Translation"
        }
    }
}

# The Status columns (Overview!E:F, zh-cn!C, de-de!C) auto-size to the
# shorter replacement text.
$newStatusColumnWidth = 12.576851254417766

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E1:F1").ColumnWidth = $newStatusColumnWidth

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C1").ColumnWidth = $newStatusColumnWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C1").ColumnWidth = $newStatusColumnWidth
